# Auto-generated: applies scheduled-runner price/profit updates to the
# Anima_Profits workbook (Sheets ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
# Each row below corresponds to one Leve entry whose market-board derived
# columns (currentAveragePrice[.NQ/.HQ], LevePrice[NQ/HQ], LeveProfit[NQ/HQ])
# were refreshed by the scheduled runner.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 3528.5715
$ws.Range("I64").Value = 3528.5715
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 3528.5715
$ws.Range("L64").Value = 0
$ws.Range("M64").Value = -3280.5715
$ws.Range("N64").ClearContents()
$ws.Range("H67").Value = 3528.5715
$ws.Range("I67").Value = 3528.5715
$ws.Range("J67").Value = 0
$ws.Range("K67").Value = 3528.5715
$ws.Range("L67").Value = 0
$ws.Range("M67").Value = -2670.5715
$ws.Range("N67").ClearContents()
$ws.Range("H101").Value = 2349.8
$ws.Range("I101").Value = 437.25
$ws.Range("J101").Value = 10000
$ws.Range("K101").Value = 1311.75
$ws.Range("L101").Value = 30000
$ws.Range("M101").Value = 310.25
$ws.Range("N101").Value = -33244
$ws.Range("H113").Value = 2661.4614
$ws.Range("I113").Value = 2457.1428
$ws.Range("J113").Value = 2899.8333
$ws.Range("K113").Value = 2457.1428
$ws.Range("L113").Value = 2899.8333
$ws.Range("M113").Value = 796.8571999999999
$ws.Range("N113").Value = -9407.8333
$ws.Range("H116").Value = 3132.6667
$ws.Range("I116").Value = 2784.1667
$ws.Range("J116").Value = 3597.3333
$ws.Range("K116").Value = 2784.1667
$ws.Range("L116").Value = 3597.3333
$ws.Range("M116").Value = 657.8332999999998
$ws.Range("N116").Value = -10481.3333
$ws.Range("H132").Value = 3068.1667
$ws.Range("I132").Value = 3098.7307
$ws.Range("K132").Value = 9296.1921
$ws.Range("M132").Value = -6766.1921
$ws.Range("H138").Value = 1928.614
$ws.Range("I138").Value = 1351.4
$ws.Range("J138").Value = 2846.9092
$ws.Range("K138").Value = 4054.2
$ws.Range("L138").Value = 8540.7276
$ws.Range("M138").Value = 1085.8
$ws.Range("N138").Value = -18820.7276
$ws.Range("H139").Value = 45500
$ws.Range("J139").Value = 46875
$ws.Range("L139").Value = 46875
$ws.Range("N139").Value = -57155
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 5116.6924
$ws.Range("I63").Value = 3835.8333
$ws.Range("J63").Value = 6214.5713
$ws.Range("K63").Value = 3835.8333
$ws.Range("L63").Value = 6214.5713
$ws.Range("M63").Value = -3149.8333
$ws.Range("N63").Value = -7586.5713
$ws.Range("H66").Value = 5116.6924
$ws.Range("I66").Value = 3835.8333
$ws.Range("J66").Value = 6214.5713
$ws.Range("K66").Value = 19179.1665
$ws.Range("L66").Value = 31072.8565
$ws.Range("M66").Value = -15747.1665
$ws.Range("N66").Value = -37936.85649999999
$ws.Range("H74").Value = 948.29266
$ws.Range("I74").Value = 641.70966
$ws.Range("K74").Value = 641.70966
$ws.Range("M74").Value = 232.29034
$ws.Range("H77").Value = 948.29266
$ws.Range("I77").Value = 641.70966
$ws.Range("K77").Value = 3208.5483
$ws.Range("M77").Value = 1159.4517
$ws.Range("H122").Value = 67902.4
$ws.Range("I122").Value = 77810.46
$ws.Range("J122").Value = 3500
$ws.Range("K122").Value = 233431.38
$ws.Range("L122").Value = 10500
$ws.Range("M122").Value = -230981.38
$ws.Range("N122").Value = -15400
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2924.4375
$ws.Range("I134").Value = 2221.2222
$ws.Range("J134").Value = 3828.5715
$ws.Range("K134").Value = 6663.6666
$ws.Range("L134").Value = 11485.7145
$ws.Range("M134").Value = -4128.6666
$ws.Range("N134").Value = -16555.7145
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5464.356
$ws.Range("I31").Value = 1229.1562
$ws.Range("J31").Value = 10483.852
$ws.Range("K31").Value = 1229.1562
$ws.Range("L31").Value = 10483.852
$ws.Range("M31").Value = -934.1561999999999
$ws.Range("N31").Value = -11073.852
$ws.Range("H34").Value = 5464.356
$ws.Range("I34").Value = 1229.1562
$ws.Range("J34").Value = 10483.852
$ws.Range("K34").Value = 1229.1562
$ws.Range("L34").Value = 10483.852
$ws.Range("M34").Value = -1027.1562
$ws.Range("N34").Value = -10887.852
$ws.Range("H122").Value = 1810.2667
$ws.Range("I122").Value = 1659.6522
$ws.Range("K122").Value = 4978.9566
$ws.Range("M122").Value = -2528.9566
$ws.Range("H123").Value = 40000
$ws.Range("J123").Value = 40000
$ws.Range("L123").Value = 40000
$ws.Range("N123").Value = -49800
$ws.Range("H134").Value = 6384.773
$ws.Range("I134").Value = 8127.857
$ws.Range("J134").Value = 3334.375
$ws.Range("K134").Value = 24383.571
$ws.Range("L134").Value = 10003.125
$ws.Range("M134").Value = -21848.571
$ws.Range("N134").Value = -15073.125
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 83333470
$ws.Range("J23").Value = 115384750
$ws.Range("L23").Value = 346154250
$ws.Range("N23").Value = -346154720
$ws.Range("H96").Value = 4600
$ws.Range("J96").Value = 4600
$ws.Range("L96").Value = 13800
$ws.Range("N96").Value = -17918
$ws.Range("H120").Value = 12166.667
$ws.Range("I120").Value = 13000
$ws.Range("K120").Value = 39000
$ws.Range("M120").Value = -34162
$ws.Range("H121").Value = 1170.1904
$ws.Range("I121").Value = 469.4
$ws.Range("J121").Value = 1264.8918
$ws.Range("K121").Value = 1408.2
$ws.Range("L121").Value = 3794.6754
$ws.Range("M121").Value = -98.19999999999982
$ws.Range("N121").Value = -6414.6754
$ws.Range("H129").Value = 2058.25
$ws.Range("I129").Value = 520
$ws.Range("K129").Value = 1560
$ws.Range("M129").Value = 3440
$ws.Range("H137").Value = 6180425
$ws.Range("J137").Value = 4626.364
$ws.Range("L137").Value = 13879.092
$ws.Range("N137").Value = -24079.092
$ws.Range("H140").Value = 1712.6666
$ws.Range("I140").Value = 1592.2222
$ws.Range("K140").Value = 4776.6666
$ws.Range("M140").Value = 403.3334000000004
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H39").Value = 19975
$ws.Range("J39").Value = 19975
$ws.Range("L39").Value = 19975
$ws.Range("N39").Value = -21039
$ws.Range("H102").Value = 1756.8572
$ws.Range("I102").Value = 1716.3334
$ws.Range("K102").Value = 1716.3334
$ws.Range("M102").Value = -94.33339999999998
$ws.Range("H122").Value = 41668444
$ws.Range("I122").Value = 1939
$ws.Range("J122").Value = 500000000
$ws.Range("K122").Value = 5817
$ws.Range("L122").Value = 1500000000
$ws.Range("M122").Value = -3367
$ws.Range("N122").Value = -1500004900
$ws.Range("H132").Value = 2583.3513
$ws.Range("I132").Value = 2097.6365
$ws.Range("J132").Value = 3295.7334
$ws.Range("K132").Value = 6292.9095
$ws.Range("L132").Value = 9887.2002
$ws.Range("M132").Value = -3762.9095
$ws.Range("N132").Value = -14947.2002
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3245
$ws.Range("I122").Value = 1816.6666
$ws.Range("J122").Value = 3857.1428
$ws.Range("K122").Value = 5449.9998
$ws.Range("L122").Value = 11571.4284
$ws.Range("M122").Value = -2999.9998
$ws.Range("N122").Value = -16471.4284
$ws.Range("H123").Value = 28666.666
$ws.Range("J123").Value = 28666.666
$ws.Range("L123").Value = 28666.666
$ws.Range("N123").Value = -38466.666
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H74").Value = 8170
$ws.Range("J74").Value = 7560
$ws.Range("L74").Value = 7560
$ws.Range("N74").Value = -9432
$ws.Range("H77").Value = 8170
$ws.Range("J77").Value = 7560
$ws.Range("L77").Value = 22680
$ws.Range("N77").Value = -32040
$ws.Range("H100").Value = 1239.05
$ws.Range("I100").Value = 1332.2142
$ws.Range("J100").Value = 1021.6667
$ws.Range("K100").Value = 2664.4284
$ws.Range("L100").Value = 2043.3334
$ws.Range("M100").Value = -2123.4284
$ws.Range("N100").Value = -3125.3334
$ws.Range("H107").Value = 582.75
$ws.Range("I107").Value = 623.1429
$ws.Range("J107").Value = 300
$ws.Range("K107").Value = 1869.4287
$ws.Range("L107").Value = 900
$ws.Range("M107").Value = 50.57129999999984
$ws.Range("N107").Value = -4740
$ws.Range("H113").Value = 903.68
$ws.Range("I113").Value = 1034.6
$ws.Range("J113").Value = 380
$ws.Range("K113").Value = 3103.8
$ws.Range("L113").Value = 1140
$ws.Range("M113").Value = -933.7999999999997
$ws.Range("N113").Value = -5480
$ws.Range("H122").Value = 1943.6207
$ws.Range("I122").Value = 1814.0385
$ws.Range("J122").Value = 3066.6667
$ws.Range("K122").Value = 5442.1155
$ws.Range("L122").Value = 9200.000100000001
$ws.Range("M122").Value = -2992.1155
$ws.Range("N122").Value = -14100.0001
$ws.Range("H126").Value = 4500
$ws.Range("I126").Value = 900
$ws.Range("K126").Value = 2700
$ws.Range("M126").Value = -230

Write-Host "Applied $(219) cell updates across ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR"
